$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (shifts current rows 4-8 down to 5-9),
# carrying formatting down from the row above (row 3), matching the
# style already used by the other data rows (e.g. date format on column D).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's data for
# "Terminal Hortofrutícola Agro Chillán" / Tuna.
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 45030
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107011
$ws.Cells.Item(4, 10).Value = "Tuna"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 40
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 18000
$ws.Cells.Item(4, 16).Value = 18000
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(4, 18).Value = "Región Metropolitana"
$ws.Cells.Item(4, 19).Value = 1000
$ws.Cells.Item(4, 20).Value = 18
